$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.172.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.803.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.288'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.64%  '
$ws.Range("E10").Value = '  +4.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0926'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.062.09'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.83%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.788.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.633'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.197.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0790'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.19%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").Value = '  -0.80%  '
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("E28").Value = '  -0.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("E32").Value = '  +1.52%  '
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.86'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.417.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.653'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.75%  '
$ws.Range("E37").Value = '  +0.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0188'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.949'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '80.67'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.24%  '
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("E43").Value = '  +3.62%  '
$ws.Range("E44").Value = '  -1.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '107.85'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0496'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.961.77'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("E48").Value = '  -1.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.39%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("E51").Value = '  +3.57%  '
